$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: "hauteur_max" (col Q) becomes a free-form text value instead of a
# plain number, and "num_siret" (col R) is re-entered as text too. ---
$ws.Range("Q2").NumberFormat = "@"
$ws.Range("Q2").Value = "290"

# --- Duplicate row 2 into a brand-new row 3 (second parking entry) ---
$ws.Range("A2:AD2").Copy()
$ws.Paste($ws.Range("A3:AD3"))

# Give the new row its own id and its own (not-applicable) hauteur_max value.
$ws.Range("A3").Value = "75114-P-002"
$ws.Range("Q3").NumberFormat = "@"
$ws.Range("Q3").Value = "N/A"

# --- Cosmetic formatting to match the "insee"/"hauteur_max"/"num_siret"
# columns being treated as text throughout (header row included). ---
$ws.Range("C1").NumberFormat = "@"
$ws.Range("Q1").NumberFormat = "@"
$ws.Range("R1").NumberFormat = "@"

# --- Turn on the header AutoFilter over the data range ---
$ws.Range("A1:AD2").AutoFilter()
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=Worksheet!`$A`$1:`$AD`$2")
$filterName.Visible = $false

# --- Leave the selection where the edit finished, like a live session would ---
$ws.Range("R3").Select()
